# "Bento object repository revisited"
#
# The FilesTab lookup query (stored as the shared string used by cell B4 on
# the "startup" sheet) is rewritten to drop the `File Type` and `Breed`
# output columns (those fields were removed from the RETURN clause, and a
# couple of trailing spaces moved around as a result of the hand-edit).
#
# Rows 2-4 of the "startup" sheet hold (in column A) the tab name and (in
# column B) the Bento query text used to populate that tab:
#   row 2 -> CasesTab
#   row 3 -> SamplesTab
#   row 4 -> FilesTab
# Only the FilesTab query text actually changes; we simply rewrite B4 with
# the new query text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newFilesTabQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["Astrocytoma"]  
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesTabQuery

# The author's selection/scroll position also moved onto that row.
[void]$ws.Range("B4").Select()
